$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Target cell values after the edit (address, new value) ---
# Column D = N for group (1); Column F = N for group (2); both recomputed for the
# corrected 2011 oxxo_counts.csv data, so every data row shifts slightly.
# Column E/G hold a handful of recomputed Mean/(SE) and Mean-difference statistics,
# and G32 (nivel_educativo_d10 mean difference) becomes significant at the 10% level.
$updates = @(
    @("D4", "16751"),
    @("E4", "2.641"),
    @("F4", "45031"),
    @("G4", "0.272***"),
    @("D6", "17621"),
    @("E6", "39.320"),
    @("F6", "45902"),
    @("G6", "1.361***"),
    @("D8", "13521"),
    @("F8", "41802"),
    @("D10", "17621"),
    @("E10", "0.450"),
    @("F10", "45902"),
    @("D12", "17621"),
    @("E12", "3.790"),
    @("F12", "45902"),
    @("G12", "-0.012"),
    @("D14", "17621"),
    @("F14", "45902"),
    @("D16", "17621"),
    @("F16", "45902"),
    @("D18", "17621"),
    @("F18", "45902"),
    @("D20", "17621"),
    @("F20", "45902"),
    @("D22", "17621"),
    @("F22", "45902"),
    @("D24", "17621"),
    @("F24", "45902"),
    @("D26", "17621"),
    @("F26", "45902"),
    @("D28", "17621"),
    @("F28", "45902"),
    @("D30", "17621"),
    @("E30", "0.200"),
    @("F30", "45902"),
    @("D32", "17621"),
    @("F32", "45902"),
    @("G32", "0.001*"),
    @("D34", "17621"),
    @("F34", "45902"),
    @("D36", "17621"),
    @("F36", "45902"),
    @("D38", "17621"),
    @("F38", "45902"),
    @("D40", "17621"),
    @("F40", "45902"),
    @("D42", "17621"),
    @("F42", "45902"),
    @("D44", "17621"),
    @("F44", "45902"),
    @("D46", "17621"),
    @("F46", "45902"),
    @("D48", "17621"),
    @("F48", "45902"),
    @("D50", "17621"),
    @("F50", "45902"),
    @("D52", "17621"),
    @("F52", "45902"),
    @("D54", "17621"),
    @("F54", "45902"),
    @("D56", "17621"),
    @("F56", "45902"),
    @("D58", "17621"),
    @("F58", "45902"),
    @("D60", "17621"),
    @("F60", "45902"),
    @("D62", "17621"),
    @("E62", "0.084"),
    @("F62", "45902"),
    @("D64", "17621"),
    @("F64", "45902"),
    @("D66", "17621"),
    @("F66", "45902"),
    @("D68", "17621"),
    @("F68", "45902"),
    @("D70", "17621"),
    @("F70", "45902"),
    @("D72", "17621"),
    @("F72", "45902"),
    @("D74", "17621"),
    @("F74", "45902"),
    @("D76", "17621"),
    @("F76", "45902"),
    @("D78", "17621"),
    @("F78", "45902"),
    @("D80", "17621"),
    @("F80", "45902"),
    @("D82", "17621"),
    @("F82", "45902"),
    @("D84", "17621"),
    @("F84", "45902"),
    @("D86", "17621"),
    @("F86", "45902"),
    @("D88", "17621"),
    @("F88", "45902"),
    @("D90", "17621"),
    @("F90", "45902"),
    @("D92", "17621"),
    @("F92", "45902"),
    @("D94", "17621"),
    @("F94", "45902"),
    @("G94", "0.027***"),
    @("D96", "17621"),
    @("F96", "45902"),
    @("D98", "17621"),
    @("F98", "45902"),
    @("D100", "17621"),
    @("F100", "45902"),
    @("D102", "17621"),
    @("F102", "45902")
)

# Values that look like plain numbers must be pinned to text format first so Excel
# keeps storing them as shared strings (matching the source workbook) instead of
# silently converting them to numeric cells.
$textRange = $ws.Range("D4,E4,F4,D6,E6,F6,D8,F8,D10,E10,F10,D12,E12,F12,G12,D14,F14,D16,F16,D18,F18,D20,F20,D22,F22,D24,F24,D26,F26,D28,F28,D30,E30,F30,D32,F32,D34,F34,D36,F36,D38,F38,D40,F40,D42,F42,D44,F44,D46,F46,D48,F48,D50,F50,D52,F52,D54,F54,D56,F56,D58,F58,D60,F60,D62,E62,F62,D64,F64,D66,F66,D68,F68,D70,F70,D72,F72,D74,F74,D76,F76,D78,F78,D80,F80,D82,F82,D84,F84,D86,F86,D88,F88,D90,F90,D92,F92,D94,F94,D96,F96,D98,F98,D100,F100,D102,F102")
foreach ($area in $textRange.Areas) {
    $area.NumberFormat = "@"
}

foreach ($update in $updates) {
    $ws.Range($update[0]).Value = $update[1]
}
